$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.532.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.71%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.612.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.91"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.76"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.76%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.88"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0998"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.332"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.76%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.083.13"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "58.455.18"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.73"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.615.04"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "333.72"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.12"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.21"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.91%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.37"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.417"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.10%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.59%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.08"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.98%  "

$ws.Range("B28").Value = "USDe"
$ws.Range("C28").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0732"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.63"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.86"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.41"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.85"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.93%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.10"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.74%  "

$ws.Range("B36").Value = "SuiNetwork"
$ws.Range("C36").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.842"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.814"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.43%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.98%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "281.61"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.81%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.594"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.04%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0941"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.04%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.94"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0527"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.86%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.941.42"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.44"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.79"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.91"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.95%  "

